$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("insertion")

# Update the raw data column A (input sizes), halving the exponent base:
# 256,1024,4096,16384,65536,262144 -> 16,32,64,128,256,512
$ws.Range("A1").Value = 16
$ws.Range("A2").Value = 32
$ws.Range("A3").Value = 64
$ws.Range("A4").Value = 128
$ws.Range("A5").Value = 256
$ws.Range("A6").Value = 512

# Update sheet view: remove the frozen/scrolled topLeftCell, change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D28").Select()

# Update the chart
$chart = $ws.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)

# Update trendline equation display (slope doubles from 1.3014 to 2.6029)
$trendline = $series.Trendlines(3)
$trendline.DisplayEquation = $true
$trendline.DataLabel.Text = "f(x) = 2,6029x - 13,546"

# Update value (x) axis minimum scale
$xAxis = $chart.Axes(1)
$xAxis.MinimumScale = 3

# Update axis titles
$xAxis.HasTitle = $true
$xAxis.AxisTitle.Text = "Taille de l'exemplaire (log22N)"

$yAxis = $chart.Axes(2)
$yAxis.HasTitle = $true
$yAxis.AxisTitle.Text = "Temps (log2y)"
